# Update the "dSF" (column F) values for a set of rows, per repulled/recalculated
# data (see commit message: "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    13 = -10
    14 = -2
    16 = -5
    17 = -4
    18 = -5
    22 = -2
    25 = 5
    30 = -1
    32 = -3
    35 = 4
    37 = 0
    40 = -3
    44 = -3
    45 = 3
    48 = 2
    55 = -1
    56 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
